$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 434.5
$ws.Range("I53").Value = 527.0909
$ws.Range("J53").Value = 230.8
$ws.Range("K53").Value = 527.0909
$ws.Range("L53").Value = 230.8
$ws.Range("M53").Value = 109.9091
$ws.Range("N53").Value = -1504.8

$ws.Range("H58").Value = 749.3333
$ws.Range("I58").Value = 186.5
$ws.Range("J58").Value = 1875
$ws.Range("K58").Value = 559.5
$ws.Range("L58").Value = 5625
$ws.Range("M58").Value = -409.5
$ws.Range("N58").Value = -5925

$ws.Range("H74").Value = 3513.3635
$ws.Range("J74").Value = 3199.1428
$ws.Range("L74").Value = 3199.1428
$ws.Range("N74").Value = -5071.1428

$ws.Range("H77").Value = 3513.3635
$ws.Range("J77").Value = 3199.1428
$ws.Range("L77").Value = 15995.714
$ws.Range("N77").Value = -25355.714

$ws.Range("H86").Value = 11120430
$ws.Range("I86").Value = 18531606
$ws.Range("J86").Value = 3668
$ws.Range("K86").Value = 18531606
$ws.Range("L86").Value = 3668
$ws.Range("M86").Value = -18530483
$ws.Range("N86").Value = -5914

$ws.Range("H89").Value = 11120430
$ws.Range("I89").Value = 18531606
$ws.Range("J89").Value = 3668
$ws.Range("K89").Value = 92658030
$ws.Range("L89").Value = 18340
$ws.Range("M89").Value = -92652414
$ws.Range("N89").Value = -29572

$ws.Range("H100").Value = 4712.222
$ws.Range("I100").Value = 3275
$ws.Range("K100").Value = 3275
$ws.Range("M100").Value = -2734

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7904.5
$ws.Range("I32").Value = 2171.1143
$ws.Range("J32").Value = 36571.43
$ws.Range("K32").Value = 2171.1143
$ws.Range("L32").Value = 36571.43
$ws.Range("M32").Value = -1884.1143
$ws.Range("N32").Value = -37145.43

$ws.Range("H61").Value = 1978.174
$ws.Range("I61").Value = 1906.3334
$ws.Range("J61").Value = 2236.8
$ws.Range("K61").Value = 1906.3334
$ws.Range("L61").Value = 2236.8
$ws.Range("M61").Value = -1694.3334
$ws.Range("N61").Value = -2660.8

$ws.Range("H97").Value = 1085.5217
$ws.Range("I97").Value = 1109.8636
$ws.Range("K97").Value = 1109.8636
$ws.Range("M97").Value = -613.8635999999999

$ws.Range("H132").Value = 1921.9302
$ws.Range("I132").Value = 1603.7428
$ws.Range("J132").Value = 3314
$ws.Range("K132").Value = 4811.2284
$ws.Range("L132").Value = 9942
$ws.Range("M132").Value = -2281.2284
$ws.Range("N132").Value = -15002

$ws.Range("H136").Value = 1978.174
$ws.Range("I136").Value = 1906.3334
$ws.Range("J136").Value = 2236.8
$ws.Range("K136").Value = 5719.0002
$ws.Range("L136").Value = 6710.400000000001
$ws.Range("M136").Value = -3169.0002
$ws.Range("N136").Value = -11810.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 50003200
$ws.Range("I86").Value = 50003200
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 50003200
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -50002077
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 50003200
$ws.Range("I89").Value = 50003200
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 250016000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -250010384
$ws.Range("N89").ClearContents()

$ws.Range("H94").Value = 13047.25
$ws.Range("I94").Value = 493.3846
$ws.Range("J94").Value = 67447.336
$ws.Range("K94").Value = 493.3846
$ws.Range("L94").Value = 67447.336
$ws.Range("M94").Value = -42.38459999999998
$ws.Range("N94").Value = -68349.336

$ws.Range("H99").Value = 1750.1666
$ws.Range("I99").Value = 1798
$ws.Range("J99").Value = 1511
$ws.Range("K99").Value = 1798
$ws.Range("L99").Value = 1511
$ws.Range("M99").Value = -300
$ws.Range("N99").Value = -4507

$ws.Range("H134").Value = 2190.423
$ws.Range("I134").Value = 1584.826
$ws.Range("J134").Value = 6833.3335
$ws.Range("K134").Value = 4754.478
$ws.Range("L134").Value = 20500.0005
$ws.Range("M134").Value = -2219.478
$ws.Range("N134").Value = -25570.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7292745
$ws.Range("I31").Value = 6484750.5
$ws.Range("J31").Value = 8336404.5
$ws.Range("K31").Value = 6484750.5
$ws.Range("L31").Value = 8336404.5
$ws.Range("M31").Value = -6484455.5
$ws.Range("N31").Value = -8336994.5

$ws.Range("H34").Value = 7292745
$ws.Range("I34").Value = 6484750.5
$ws.Range("J34").Value = 8336404.5
$ws.Range("K34").Value = 6484750.5
$ws.Range("L34").Value = 8336404.5
$ws.Range("M34").Value = -6484548.5
$ws.Range("N34").Value = -8336808.5

$ws.Range("H58").Value = 1891.1177
$ws.Range("I58").Value = 1036.9
$ws.Range("J58").Value = 3111.4285
$ws.Range("K58").Value = 1036.9
$ws.Range("L58").Value = 3111.4285
$ws.Range("M58").Value = -833.9000000000001
$ws.Range("N58").Value = -3517.4285

$ws.Range("H136").Value = 1891.1177
$ws.Range("I136").Value = 1036.9
$ws.Range("J136").Value = 3111.4285
$ws.Range("K136").Value = 3110.7
$ws.Range("L136").Value = 9334.2855
$ws.Range("M136").Value = -560.7000000000003
$ws.Range("N136").Value = -14434.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 34486588
$ws.Range("I64").Value = 3900
$ws.Range("J64").Value = 35718110
$ws.Range("K64").Value = 11700
$ws.Range("L64").Value = 107154330
$ws.Range("M64").Value = -11430
$ws.Range("N64").Value = -107154870

$ws.Range("H67").Value = 34486588
$ws.Range("I67").Value = 3900
$ws.Range("J67").Value = 35718110
$ws.Range("K67").Value = 11700
$ws.Range("L67").Value = 107154330
$ws.Range("M67").Value = -10764
$ws.Range("N67").Value = -107156202

$ws.Range("H106").Value = 7583.7393
$ws.Range("J106").Value = 7821.3
$ws.Range("L106").Value = 23463.9
$ws.Range("N106").Value = -25355.9

$ws.Range("H131").Value = 994.2364
$ws.Range("I131").Value = 739.1818
$ws.Range("J131").Value = 1058
$ws.Range("K131").Value = 2217.5454
$ws.Range("L131").Value = 3174
$ws.Range("M131").Value = 2822.4546
$ws.Range("N131").Value = -13254

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1538.625
$ws.Range("I97").Value = 1538.625
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1538.625
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1042.625
$ws.Range("N97").ClearContents()

$ws.Range("H132").Value = 2436.9473
$ws.Range("I132").Value = 1686.9333
$ws.Range("J132").Value = 5249.5
$ws.Range("K132").Value = 5060.7999
$ws.Range("L132").Value = 15748.5
$ws.Range("M132").Value = -2530.7999
$ws.Range("N132").Value = -20808.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2169.3
$ws.Range("I93").Value = 1961.625
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 1961.625
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -713.625
$ws.Range("N93").Value = -5496

$ws.Range("H132").Value = 1232509.9
$ws.Range("I132").Value = 1961954.8
$ws.Range("J132").Value = 3971.2104
$ws.Range("K132").Value = 5885864.4
$ws.Range("L132").Value = 11913.6312
$ws.Range("M132").Value = -5883334.4
$ws.Range("N132").Value = -16973.6312

$ws.Range("H136").Value = 4037237
$ws.Range("I136").Value = 5440542
$ws.Range("J136").Value = 2735.625
$ws.Range("K136").Value = 16321626
$ws.Range("L136").Value = 8206.875
$ws.Range("M136").Value = -16319076
$ws.Range("N136").Value = -13306.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 90911550
$ws.Range("I81").Value = 333334600
$ws.Range("J81").Value = 2912
$ws.Range("K81").Value = 666669200
$ws.Range("L81").Value = 5824
$ws.Range("M81").Value = -666668139
$ws.Range("N81").Value = -7946

$ws.Range("H84").Value = 90911550
$ws.Range("I84").Value = 333334600
$ws.Range("J84").Value = 2912
$ws.Range("K84").Value = 3333346000
$ws.Range("L84").Value = 29120
$ws.Range("M84").Value = -3333340696
$ws.Range("N84").Value = -39728

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H96").Value = 1001.5
$ws.Range("I96").Value = 1002
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 1002
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = 371
$ws.Range("N96").Value = -3746

$ws.Range("H132").Value = 3737
$ws.Range("I132").Value = 3015
$ws.Range("J132").Value = 5301.3335
$ws.Range("K132").Value = 9045
$ws.Range("L132").Value = 15904.0005
$ws.Range("M132").Value = -6515
$ws.Range("N132").Value = -20964.0005

